# Economic Dashboard weekly data refresh - 2025-12-11
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- GDP Nowcast (F7) ---
$ws.Range("F7").Value = 0.215816678152998

# --- Row 13: UI Initial Claims (ICSA) ---
$ws.Range("N13").Value = 45992
$ws.Range("Q13").Value = 236000
$ws.Range("R13").Value = 192000
$ws.Range("S13").Value = 217000
$ws.Range("T13").Value = 222000
$ws.Range("U13").Value = 228000

# --- Row 14: UI Continuing Claims (CCSA) ---
$ws.Range("N14").Value = 45985
$ws.Range("Q14").Value = 1838000
$ws.Range("R14").Value = 1937000
$ws.Range("S14").Value = 1944000
$ws.Range("T14").Value = 1953000
$ws.Range("U14").Value = 1946000

# --- Clear the "new this week" highlight on C28/C29/C30/C31/N51 ---
# (these date cells move from the highlighted style (s=48) back to the
# plain style (s=47); copy formats from a cell that already carries s=47)
$ws.Range("N7").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("N51").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 29: 5yr, 5yr Forward (T5YIFR) ---
$ws.Range("N29").Value = 46001
$ws.Range("Q29").Value = 2.18
$ws.Range("S29").Value = 2.2
$ws.Range("T29").Value = $null
$ws.Range("U29").Value = $null

# --- Row 30: 10yr TIPS (T10YIE) ---
$ws.Range("N30").Value = 46001
$ws.Range("Q30").Value = 2.25
$ws.Range("S30").Value = 2.26
$ws.Range("T30").Value = $null
$ws.Range("U30").Value = $null

# --- N47: FFR (DFF) date refresh only ---
$ws.Range("N47").Value = 46000

# --- Row 48: 2y UST (DGS2) ---
$ws.Range("N48").Value = 46000
$ws.Range("Q48").Value = 3.61
$ws.Range("R48").Value = 3.57
$ws.Range("S48").Value = $null
$ws.Range("T48").Value = $null
$ws.Range("U48").Value = 3.56

# --- Row 49: 5y UST (DGS5) ---
$ws.Range("N49").Value = 46000
$ws.Range("Q49").Value = 3.78
$ws.Range("R49").Value = 3.75
$ws.Range("S49").Value = $null
$ws.Range("T49").Value = $null
$ws.Range("U49").Value = 3.72

# --- Row 50: 10y UST (DGS10) ---
$ws.Range("N50").Value = 46000
$ws.Range("Q50").Value = 4.18
$ws.Range("R50").Value = 4.17
$ws.Range("S50").Value = $null
$ws.Range("T50").Value = $null
$ws.Range("U50").Value = 4.14

# --- Row 52: BAA (DBAA) ---
$ws.Range("N52").Value = 46000
$ws.Range("Q52").Value = 5.91
$ws.Range("R52").Value = 5.9
$ws.Range("S52").Value = $null
$ws.Range("T52").Value = $null
$ws.Range("U52").Value = 5.88
